$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g3 = $s.Shapes.Item(1)
$zero = $g3.GroupItems.Item("Rounded Rectangle 48")
$tr = $zero.TextFrame.TextRange
$tr.Text = "Kinesis"
$tr.Font.Bold = $true
Write-Host "done"
